# Add a second worksheet ("Sheet2") after "Sheet1" for
# "Create Standard Parts From Excel", populate its header row and six
# sample data rows, size its columns, and leave it as the active sheet
# (mirrors a user inserting + filling in a new tab in the Excel UI).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert right after Sheet1 so sheet order is Sheet1, Sheet2 and the
# freshly-added sheet becomes the active tab.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Header row.
$headers = @("DRG NO.", "Line1", "Line2", "Character Height", "Dim A", "Dim B", "Remarks")
for ($c = 1; $c -le $headers.Length; $c++) {
    $ws2.Cells.Item(1, $c).Value = $headers[$c - 1]
}

# Sample data - six rows. Shared strings must land in "column major"
# order (all of column A, then all of column B) to match how the
# workbook was authored.
$drgNumbers = @("DT-123456-001", "DT-123456-002", "DT-123456-003", "DT-123456-004", "DT-123456-005", "DT-123456-006")
$lines = @("LOREM", "IPSUM", "DOLOR", "SIT", "AMET", "CONSECTETUR")

for ($r = 0; $r -lt $drgNumbers.Length; $r++) {
    $ws2.Cells.Item($r + 2, 1).Value = $drgNumbers[$r]
}
for ($r = 0; $r -lt $lines.Length; $r++) {
    $ws2.Cells.Item($r + 2, 2).Value = $lines[$r]
}
for ($r = 0; $r -lt $drgNumbers.Length; $r++) {
    $row = $r + 2
    $ws2.Cells.Item($row, 4).Value = 0.7
    $ws2.Cells.Item($row, 5).Value = 50
    $ws2.Cells.Item($row, 6).Value = 25
}

# Column widths (best-fit-sized by the original author).
$ws2.Columns.Item(1).ColumnWidth = 13.022135416666666
$ws2.Columns.Item(2).ColumnWidth = 12.877604166666666
$ws2.Columns.Item(3).ColumnWidth = 4.877604166666667
$ws2.Columns.Item(4).ColumnWidth = 15.022135416666666
$ws2.Columns.Item(5).ColumnWidth = 5.451822916666667
$ws2.Columns.Item(6).ColumnWidth = 5.307291666666667
$ws2.Columns.Item(7).ColumnWidth = 7.736979166666667
